# Weekly refresh of the NFL predictions sheet ("Tabelle1"):
#  - New week's four games (teams + Vegas spread/total) in B2:D5
#  - The two already-played games' final scores (J3, J4, J5) and the
#    new team label B4 (DET/MIN week -> MIN week)
#  - The rolling 11-game "Winner" (L2:L12) / "Loser" (S2:S12) score
#    history used by the SMALL/MATCH/INDEX lookups is shifted to add
#    this week's new result and drop the oldest one.
# Everything else (G1/J1 totals, E:I/M:X helper columns, B9:C12
# summaries, row 16-28 frequency table) is formula-driven and
# recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: NE vs BUF ---
$ws.Range("B2").Value = "BUF"
$ws.Range("C2").Value = -13
$ws.Range("D2").Value = 40.5

# --- Row 3: KC game (opponent label unchanged) ---
$ws.Range("C3").Value = -7
$ws.Range("D3").Value = 44.5
$ws.Range("J3").Value = 20

# --- Row 4: now MIN (was DET) ---
$ws.Range("B4").Value = "MIN"
$ws.Range("C4").Value = -2
$ws.Range("D4").Value = 46.5
$ws.Range("J4").Value = 22

# --- Row 5: LAR game (opponent label unchanged) ---
$ws.Range("C5").Value = -5.5
$ws.Range("D5").Value = 44.5
$ws.Range("J5").Value = 20

# --- L2:L12 "Winner" score-pair history (rolling window) ---
$ws.Range("L2").Value  = "  27   67"
$ws.Range("L3").Value  = "  20   56"
$ws.Range("L4").Value  = "  24   54"
$ws.Range("L5").Value  = "  31   49"
$ws.Range("L6").Value  = "  30   43"
$ws.Range("L7").Value  = "  34   36"
$ws.Range("L8").Value  = "  23   35"
$ws.Range("L9").Value  = "  28   28"
$ws.Range("L10").Value = "  17   27"
$ws.Range("L11").Value = "  37   26"
$ws.Range("L12").Value = "  19   24"

# --- S2:S12 "Loser" score-pair history (rolling window) ---
$ws.Range("S2").Value  = "  17   80"
$ws.Range("S3").Value  = "  10   78"
$ws.Range("S4").Value  = "  20   52"
$ws.Range("S5").Value  = "  16   49"
$ws.Range("S6").Value  = "  13   48"
$ws.Range("S7").Value  = "  14   34"
$ws.Range("S8").Value  = "  24   33"
$ws.Range("S9").Value  = "  3    30"
$ws.Range("S10").Value = "  22   27"
$ws.Range("S11").Value = "  21   25"
$ws.Range("S12").Value = "  6    24"

Write-Host "Sheet updated"
